$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Enter the roster (names + scores + averages) first -- this mirrors the
#    authoring order implied by the shared-string table (names before
#    headers).
# ---------------------------------------------------------------------------

# Row 3 - Marvin
$ws.Range("C3").Value = "Marvin"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 20
$ws.Range("F3").Formula = "=+AVERAGE(D3:E3)"

# Row 4 - Jesús
$ws.Range("C4").Value = "Jesús"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 17
$ws.Range("F4").Formula = "=+AVERAGE(D4:E4)"

# Row 5 - Jenny
$ws.Range("C5").Value = "Jenny"
$ws.Range("D5").Value = 18
$ws.Range("E5").Value = 16
$ws.Range("F5").Formula = "=+AVERAGE(D5:E5)"

# Row 6 - Carol
$ws.Range("C6").Value = "Carol"
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 16
$ws.Range("F6").Formula = "=+AVERAGE(D6:E6)"

# Row 7 - Iris
$ws.Range("C7").Value = "Iris"
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = 17
$ws.Range("F7").Formula = "=+AVERAGE(D7:E7)"

# Row 8 - Kevin
$ws.Range("C8").Value = "Kevin"
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 16
$ws.Range("F8").Formula = "=+AVERAGE(D8:E8)"

# Row 9 - Jose Uribe
$ws.Range("C9").Value = "Jose Uribe"
$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 17
$ws.Range("F9").Formula = "=+AVERAGE(D9:E9)"

# ---------------------------------------------------------------------------
# 2) Now add the header row, in the column order C, E, F, D (this is the
#    order that reproduces the shared-string indices seen in the target).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Representante de Grupo"
$ws.Range("E2").Value = "Exposición"
$ws.Range("F2").Value = "Nota PC2"
$ws.Range("D2").Value = "Trabajo"

# ---------------------------------------------------------------------------
# 3) Styling -- borders, bold names, center/middle alignment, header wrap.
# ---------------------------------------------------------------------------

# Bold the name column (C3:C9) and give every data/header cell a thin box
# border.
$ws.Range("C3:C9").Font.Bold = $true
$ws.Range("C2:F9").Borders.LineStyle = 1

# Row 3 D:E -- center + middle
$ws.Range("D3:E3").HorizontalAlignment = -4108
$ws.Range("D3:E3").VerticalAlignment = -4108

# Row 4:5 D:E -- center only (horizontal), matches the author's slightly
# different formatting pass on those two rows.
$ws.Range("D4:E5").HorizontalAlignment = -4108

# Row 6:9 D:E -- center + middle again
$ws.Range("D6:E9").HorizontalAlignment = -4108
$ws.Range("D6:E9").VerticalAlignment = -4108

# Averages column (F) -- bold, center + middle
$ws.Range("F3:F9").Font.Bold = $true
$ws.Range("F3:F9").HorizontalAlignment = -4108
$ws.Range("F3:F9").VerticalAlignment = -4108

# Header row -- bold, center + middle
$ws.Range("C2:F2").Font.Bold = $true
$ws.Range("C2:F2").HorizontalAlignment = -4108
$ws.Range("C2:F2").VerticalAlignment = -4108
$ws.Range("C2").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Row height / column widths / selection.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 30

$ws.Columns("A").ColumnWidth = 2.1666666
$ws.Columns("C").ColumnWidth = 13.1666666
$ws.Columns("D:E").ColumnWidth = 10.4166666

$ws.Range("C12").Select()
